# Natmi following Dr Hou advice
# Update LR-pair data rows (sheet1) to new sending/target cluster breakdown including FAPs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Lgi1"
$ws.Range("C2").Value = "Adam22"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3953366666666667
$ws.Range("H2").Value = 1.18601
$ws.Range("I2").Value = 0.5534398610532347
$ws.Range("J2").Value = 0.5534398610532348
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8328803333333332
$ws.Range("N2").Value = 2.498641
$ws.Range("O2").Value = 0.03917234793046614
$ws.Range("P2").Value = 0.03917234793046613
$ws.Range("Q2").Value = 0.3292681347122222
$ws.Range("R2").Value = 2.963413212409999
$ws.Range("S2").Value = 0.02167953879576615
$ws.Range("T2").Value = 0.02167953879576615

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Lgi1"
$ws.Range("C3").Value = "Adam22"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3953366666666667
$ws.Range("H3").Value = 1.18601
$ws.Range("I3").Value = 0.5534398610532347
$ws.Range("J3").Value = 0.5534398610532348
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.873409
$ws.Range("N3").Value = 14.620227
$ws.Range("O3").Value = 0.2292080450398417
$ws.Range("P3").Value = 0.2292080450398417
$ws.Range("Q3").Value = 1.926637269363333
$ws.Range("R3").Value = 17.33973542427
$ws.Range("S3").Value = 0.1268528685991336
$ws.Range("T3").Value = 0.1268528685991336

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Lgi1"
$ws.Range("C4").Value = "Adam22"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3953366666666667
$ws.Range("H4").Value = 1.18601
$ws.Range("I4").Value = 0.5534398610532347
$ws.Range("J4").Value = 0.5534398610532348
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 15.55565633333333
$ws.Range("N4").Value = 46.66696899999999
$ws.Range("O4").Value = 0.7316196070296923
$ws.Range("P4").Value = 0.7316196070296922
$ws.Range("Q4").Value = 6.149721322632222
$ws.Range("R4").Value = 55.34749190368999
$ws.Range("S4").Value = 0.4049074536583351
$ws.Range("T4").Value = 0.4049074536583351

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Lgi1"
$ws.Range("C5").Value = "Adam22"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3189896666666667
$ws.Range("H5").Value = 0.956969
$ws.Range("I5").Value = 0.4465601389467652
$ws.Range("J5").Value = 0.4465601389467652
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8328803333333332
$ws.Range("N5").Value = 2.498641
$ws.Range("O5").Value = 0.03917234793046614
$ws.Range("P5").Value = 0.03917234793046613
$ws.Range("Q5").Value = 0.2656802199032222
$ws.Range("R5").Value = 2.391121979129
$ws.Range("S5").Value = 0.01749280913469999
$ws.Range("T5").Value = 0.01749280913469999

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Lgi1"
$ws.Range("C6").Value = "Adam22"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3189896666666667
$ws.Range("H6").Value = 0.956969
$ws.Range("I6").Value = 0.4465601389467652
$ws.Range("J6").Value = 0.4465601389467652
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.873409
$ws.Range("N6").Value = 14.620227
$ws.Range("O6").Value = 0.2292080450398417
$ws.Range("P6").Value = 0.2292080450398417
$ws.Range("Q6").Value = 1.554567112440333
$ws.Range("R6").Value = 13.991104011963
$ws.Range("S6").Value = 0.1023551764407081
$ws.Range("T6").Value = 0.1023551764407081

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Lgi1"
$ws.Range("C7").Value = "Adam22"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3189896666666667
$ws.Range("H7").Value = 0.956969
$ws.Range("I7").Value = 0.4465601389467652
$ws.Range("J7").Value = 0.4465601389467652
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 15.55565633333333
$ws.Range("N7").Value = 46.66696899999999
$ws.Range("O7").Value = 0.7316196070296923
$ws.Range("P7").Value = 0.7316196070296922
$ws.Range("Q7").Value = 4.962093628551222
$ws.Range("R7").Value = 44.65884265696099
$ws.Range("S7").Value = 0.3267121533713572
$ws.Range("T7").Value = 0.3267121533713571

